# MBD Problem Solver Input.xlsx - update test case for singular axial torque
# causing torsion (both angle of twist and torsional shear stress work)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cross-section dimension (L2): 0.0125 -> 0.015
$ws.Range("L2").Value = 0.015

# Axial torque magnitude (O2): 3.18 -> 50
$ws.Range("O2").Value = 50

# Torsion variable - target radius from neutral axis (Q2): 1 -> 0.015
$ws.Range("Q2").Value = 0.015

# Row 3: C-S dimension second value (L3): 0 -> blank
$ws.Range("L3").ClearContents()

# Row 3: Axial torque magnitude (O3): 6.37 -> blank
$ws.Range("O3").ClearContents()

# Row 3: Position of AT (P3): 1 -> blank
$ws.Range("P3").ClearContents()

# Row 3: Torsion variable - total length of beam (Q3): 2 -> 1
$ws.Range("Q3").Value = 1

# Row 3: Moduli - Shear modulus, now a formula (R3): 20000 -> =27*10^9
$ws.Range("R3").Formula = "=27*10^9"

# Row 4: Axial torque magnitude (O4): 9.55 -> blank
$ws.Range("O4").ClearContents()

# Row 4: Position of AT (P4): 2 -> blank
$ws.Range("P4").ClearContents()

# Update shared-string text for the Moduli description header (R10)
$ws.Range("R10").Value = "First cell: Young's Modulus value, E`nSecond cell: Shear Modulus value, G"

# Update selection to match final cursor position
$ws.Range("M3").Select() | Out-Null
